# Updated cryptos list on Wed Sep 25 04:33:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'64.364.85"
$ws.Range("E2").Value = "`'  +2.01%  "
$ws.Range("D3").Value = "`'2.634.97"
$ws.Range("E3").Value = "`'  +0.28%  "
$ws.Range("E4").Value = "`'  +0.09%  "
$ws.Range("D5").Value = "`'605.41"
$ws.Range("E5").Value = "`'  +0.08%  "
$ws.Range("D6").Value = "`'151.33"
$ws.Range("E6").Value = "`'  +3.47%  "
$ws.Range("E7").Value = "`'  +0.03%  "
$ws.Range("D8").Value = "`'0.592"
$ws.Range("E8").Value = "`'  +1.14%  "
$ws.Range("D9").Value = "`'0.110"
$ws.Range("E9").Value = "`'  +2.40%  "
$ws.Range("D10").Value = "`'5.78"
$ws.Range("E10").Value = "`'  +3.29%  "
$ws.Range("D11").Value = "`'0.386"
$ws.Range("E11").Value = "`'  +6.55%  "
$ws.Range("E12").Value = "`'  -0.62%  "
$ws.Range("D13").Value = "`'27.75"
$ws.Range("E13").Value = "`'  +2.13%  "
$ws.Range("D14").Value = "`'3.111.74"
$ws.Range("E14").Value = "`'  +0.58%  "
$ws.Range("D15").Value = "`'64.225.15"
$ws.Range("E15").Value = "`'  +2.03%  "
$ws.Range("E16").Value = "`'  +4.37%  "
$ws.Range("D17").Value = "`'2.641.49"
$ws.Range("E17").Value = "`'  +1.00%  "
$ws.Range("D18").Value = "`'12.19"
$ws.Range("E18").Value = "`'  +8.05%  "
$ws.Range("D19").Value = "`'4.66"
$ws.Range("E19").Value = "`'  +4.42%  "
$ws.Range("D20").Value = "`'350.63"
$ws.Range("E20").Value = "`'  +3.38%  "
$ws.Range("D21").Value = "`'7.01"
$ws.Range("E21").Value = "`'  +2.19%  "
$ws.Range("E22").Value = "`'  +0.33%  "
$ws.Range("D23").Value = "`'5.73"
$ws.Range("E23").Value = "`'  +3.06%  "
$ws.Range("D24").Value = "`'66.81"
$ws.Range("E24").Value = "`'  +0.50%  "
$ws.Range("D25").Value = "`'1.75"
$ws.Range("E25").Value = "`'  +14.47%  "
$ws.Range("D26").Value = "`'1.71"
$ws.Range("E26").Value = "`'  +5.48%  "
$ws.Range("D27").Value = "`'9.31"
$ws.Range("E27").Value = "`'  +7.52%  "
$ws.Range("B28").Value = "`'Kaspa"
$ws.Range("C28").Value = "`'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "`'0.167"
$ws.Range("E28").Value = "`'  +2.66%  "
$ws.Range("B29").Value = "`'Aptos"
$ws.Range("C29").Value = "`'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "`'8.17"
$ws.Range("E29").Value = "`'  +3.48%  "
$ws.Range("D30").Value = "`'547.39"
$ws.Range("E30").Value = "`'  +1.98%  "
$ws.Range("E31").Value = "`'  +0.04%  "
$ws.Range("E32").Value = "`'  +2.47%  "
$ws.Range("D33").Value = "`'0.0₃0864"
$ws.Range("E33").Value = "`'  +7.97%  "
$ws.Range("D34").Value = "`'1.77"
$ws.Range("E34").Value = "`'  +1.15%  "
$ws.Range("D35").Value = "`'5.30"
$ws.Range("E35").Value = "`'  +1.09%  "
$ws.Range("D36").Value = "`'167.52"
$ws.Range("E36").Value = "`'  -1.00%  "
$ws.Range("E37").Value = "`'  +8.68%  "
$ws.Range("D38").Value = "`'0.412"
$ws.Range("E38").Value = "`'  +2.30%  "
$ws.Range("E39").Value = "`'  -0.02%  "
$ws.Range("D40").Value = "`'19.60"
$ws.Range("E40").Value = "`'  +3.22%  "
$ws.Range("B41").Value = "`'Aave"
$ws.Range("C41").Value = "`'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "`'171.64"
$ws.Range("E41").Value = "`'  +2.28%  "
$ws.Range("B42").Value = "`'USDe"
$ws.Range("C42").Value = "`'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "`'1.00"
$ws.Range("E42").Value = "`'  +0.06%  "
$ws.Range("D43").Value = "`'40.06"
$ws.Range("E43").Value = "`'  +0.90%  "
$ws.Range("D44").Value = "`'3.95"
$ws.Range("E44").Value = "`'  +5.74%  "
$ws.Range("D45").Value = "`'0.0588"
$ws.Range("E45").Value = "`'  +4.07%  "
$ws.Range("D46").Value = "`'21.65"
$ws.Range("E46").Value = "`'  -2.54%  "
$ws.Range("D47").Value = "`'0.630"
$ws.Range("E47").Value = "`'  +1.06%  "
$ws.Range("D48").Value = "`'2.01"
$ws.Range("E48").Value = "`'  +14.56%  "
$ws.Range("D49").Value = "`'0.0246"
$ws.Range("E49").Value = "`'  +2.63%  "
$ws.Range("D50").Value = "`'0.0967"
$ws.Range("E50").Value = "`'  +1.01%  "
$ws.Range("D51").Value = "`'19.38"
$ws.Range("E51").Value = "`'  +4.63%  "
